$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.737.84'
$ws.Range('E2').Value = '  +15.01%  '

$ws.Range('D3').Value = '1.727.09'
$ws.Range('E3').Value = '  +8.54%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9940'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.45%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '312.29'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.52%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9882'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.18%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3757'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +4.30%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '50.53'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +23.24%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3553'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +6.76%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.213'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +8.87%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07481'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +8.32%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.9892'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.68%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.59'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +11.80%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.312'
$ws.Range('D14').ClearFormats()

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.923'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +6.78%  '

$ws.Range('D16').Value = '1.724.28'
$ws.Range('E16').Value = '  +8.32%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001136'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +7.24%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.9879'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.27%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06737'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.54%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '85.52'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +12.18%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.12'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +8.43%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.292'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +7.06%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.45'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +8.77%  '

$ws.Range('D24').Value = '25.634.42'
$ws.Range('E24').Value = '  +14.54%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.420'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.22%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.824'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +12.80%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '154.18'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.73%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.24'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +6.32%  '

$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '132.66'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +8.29%  '

$ws.Range('B30').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C30').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D30').Value = '1.921.84'
$ws.Range('E30').Value = '  +9.01%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.149'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +25.24%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.600'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +12.64%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.081'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.62%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.770'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +8.31%  '

$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '13.23'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +13.32%  '

$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08510'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +5.22%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06590'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +10.29%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.507'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +8.51%  '

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '9.058'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +8.66%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.02409'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +10.71%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.2164'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +10.06%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.267'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.98%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.6384'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +11.06%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.9878'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.13%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.54'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +7.15%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6168'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +10.59%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.850'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.61%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.113'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +9.60%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '130.54'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.63%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07431'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +9.99%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '78.00'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +8.15%  '
